$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$old = [char]0x201C + "2022. Año del Quincentenario de Toluca, Capital del Estado de México" + [char]0x22 + "."
$new = '${leyenda}'

$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
